$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 (17:22 -> 17:52)
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 17:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1068562
$ws.Range("C4").Value = 4368
$ws.Range("D4").Value = 148013
$ws.Range("E4").Value = 858374
$ws.Range("G4").Value = 519
$ws.Range("H4").Value = 62175

# --- Row 9: Alemania ---
$ws.Range("B9").Value = 162123
$ws.Range("C9").Value = 584
$ws.Range("E9").Value = 32105
$ws.Range("G9").Value = 51
$ws.Range("H9").Value = 6518

# --- Row 35: Polonia ---
$ws.Range("B35").Value = 12877
$ws.Range("C35").Value = 237
$ws.Range("E35").Value = 8997
$ws.Range("G35").Value = 20
$ws.Range("H35").Value = 644

# --- Row 46: Chequia ---
$ws.Range("B46").Value = 7642
$ws.Range("C46").Value = 63
$ws.Range("D46").Value = 3281
$ws.Range("E46").Value = 4126
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 235

# --- Rows 52/53: Sudafrica & Egipto swap order with updated Egipto numbers ---
# Row 52 used to be Sudafrica, row 53 used to be Egipto.
# After the edit, Egipto (with refreshed figures) sorts ahead of Sudafrica,
# so row 52 becomes Egipto (new data) and row 53 becomes Sudafrica
# (keeping its previous, unchanged figures).
$ws.Range("A52").Value = "Egipto"
$ws.Range("B52").Value = 5537
$ws.Range("C52").Value = 269
$ws.Range("D52").Value = 1381
$ws.Range("E52").Value = 3764
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 12
$ws.Range("H52").Value = 392

$ws.Range("A53").Value = "Sudafrica"
$ws.Range("B53").Value = 5350
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 2073
$ws.Range("E53").Value = 3174
$ws.Range("F53").Value = 36
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 103

# --- Row 115: Jordania ---
$ws.Range("B115").Value = 453
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 362
$ws.Range("E115").Value = 83

# --- Row 126: Montenegro ---
$ws.Range("D126").Value = 214
$ws.Range("E126").Value = 101

# --- Row 128: Isla de Man ---
$ws.Range("D128").Value = 260
$ws.Range("E128").Value = 34
